# 116Cd(d,p)117Cd spectroscopic-factor table: refreshed fit results.
# Existing data rows (sheet rows 2-34, A:E) get new values, and a new
# 34th data point (sheet row 35) is appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2:34 (A:E) with new data ---
$data = New-Object 'object[,]' 33,5
$data[0,0] = 0
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0.3883432743546536
$data[0,4] = 0.005801046176978834
$data[1,0] = 1
$data[1,1] = 107.8585558484507
$data[1,2] = 2
$data[1,3] = 0.009496397978399362
$data[1,4] = 0
$data[2,0] = 2
$data[2,1] = 107.8585558484507
$data[2,2] = 5
$data[2,3] = 0.001601448494939445
$data[2,4] = 0
$data[3,0] = 3
$data[3,1] = 134.5163962667762
$data[3,2] = 2
$data[3,3] = 0.4223253870688833
$data[3,4] = 0
$data[4,0] = 4
$data[4,1] = 134.5163962667762
$data[4,2] = 5
$data[4,3] = 0.03594088430122813
$data[4,4] = 0
$data[5,0] = 5
$data[5,1] = 292.2372701306018
$data[5,2] = 0
$data[5,3] = 0.01994287695521927
$data[5,4] = 0.001313111239849829
$data[6,0] = 6
$data[6,1] = 411.090801103067
$data[6,2] = 2
$data[6,3] = 0.00224922927944456
$data[6,4] = 0.0002320633383553911
$data[7,0] = 7
$data[7,1] = 426.2
$data[7,2] = 0
$data[7,3] = 0.06328771551655087
$data[7,4] = 0.002759137461212026
$data[8,0] = 8
$data[8,1] = 442.6
$data[8,2] = 2
$data[8,3] = 0.04028456662544167
$data[8,4] = 0.0008912514740141964
$data[9,0] = 9
$data[9,1] = 498.0921546689657
$data[9,2] = 5
$data[9,3] = 0.1286545273913219
$data[9,4] = 0.003116888034738211
$data[10,0] = 10
$data[10,1] = 522.1
$data[10,2] = 2
$data[10,3] = 0.02843631200254661
$data[10,4] = 0.000734517284198621
$data[11,0] = 11
$data[11,1] = 627.2906398511834
$data[11,2] = 5
$data[11,3] = 0.01233653079092983
$data[11,4] = 0.00110379486024109
$data[12,0] = 12
$data[12,1] = 639.4806312796794
$data[12,2] = 5
$data[12,3] = 0.008879517077967592
$data[12,4] = 0.0009722098990475464
$data[13,0] = 13
$data[13,1] = 665.2
$data[13,2] = 2
$data[13,3] = 0.06399918120629153
$data[13,4] = 0.001124286920312814
$data[14,0] = 14
$data[14,1] = 690.8
$data[14,2] = 2
$data[14,3] = 0.01065132706000522
$data[14,4] = 0.000576664203885633
$data[15,0] = 15
$data[15,1] = 699.9918964602373
$data[15,2] = 0
$data[15,3] = 0.003474240466262073
$data[15,4] = 0.001789760240195614
$data[16,0] = 16
$data[16,1] = 728.2049555816418
$data[16,2] = 2
$data[16,3] = 0.0007583871226285256
$data[16,4] = 0.0001516774245257051
$data[17,0] = 17
$data[17,1] = 779.0788707099858
$data[17,2] = 0
$data[17,3] = 0.01325238894025671
$data[17,4] = 0.001356543749790057
$data[18,0] = 18
$data[18,1] = 820.1
$data[18,2] = 2
$data[18,3] = 0.01292736193670261
$data[18,4] = 0.0004812496741519585
$data[19,0] = 19
$data[19,1] = 864.8217727435604
$data[19,2] = 5
$data[19,3] = 0.006112995676851511
$data[19,4] = 0.0007485300828797769
$data[20,0] = 20
$data[20,1] = 997.3547572460591
$data[20,2] = 5
$data[20,3] = 0.005240126451864466
$data[20,4] = 0.0007921121380725358
$data[21,0] = 21
$data[21,1] = 1011.995701247272
$data[21,2] = 3
$data[21,3] = 0.001929248399486033
$data[21,4] = 0.0004823120998715085
$data[22,0] = 22
$data[22,1] = 1073.2
$data[22,2] = 0
$data[22,3] = 0.01135728761882676
$data[22,4] = 0.001323179140057487
$data[23,0] = 23
$data[23,1] = 1082.466210769789
$data[23,2] = 5
$data[23,3] = 0.02441890964009653
$data[23,4] = 0.001559930345558992
$data[24,0] = 24
$data[24,1] = 1105.246920920629
$data[24,2] = 1
$data[24,3] = 0.005243685527966736
$data[24,4] = 0.0003343261806634169
$data[25,0] = 25
$data[25,1] = 1132.989830446533
$data[25,2] = 3
$data[25,3] = 0.009050998742610073
$data[25,4] = 0.0009120918663770896
$data[26,0] = 26
$data[26,1] = 1221.150820773939
$data[26,2] = 0
$data[26,3] = 0.02144127001002603
$data[26,4] = 0.00175241149120405
$data[27,0] = 27
$data[27,1] = 1257.056227090943
$data[27,2] = 0
$data[27,3] = 0.01604363934744349
$data[27,4] = 0.001542657629561874
$data[28,0] = 28
$data[28,1] = 1278.325416780053
$data[28,2] = 2
$data[28,3] = 0.01383871555435027
$data[28,4] = 0.0004898660373221334
$data[29,0] = 29
$data[29,1] = 1315.93448102085
$data[29,2] = 0
$data[29,3] = 0.00674385453734916
$data[29,4] = 0.0011055499241556
$data[30,0] = 30
$data[30,1] = 1342.277962181449
$data[30,2] = 2
$data[30,3] = 0.003600336193038126
$data[30,4] = 0.0003176767229151288
$data[31,0] = 31
$data[31,1] = 1355.9
$data[31,2] = 2
$data[31,3] = 0.01311240045682935
$data[31,4] = 0.0004979392578542791
$data[32,0] = 32
$data[32,1] = 1475.089127951248
$data[32,2] = 2
$data[32,3] = 0.01791810606136577
$data[32,4] = 0.0006499560319044467
$ws.Range("A2:E34").Value = $data

# --- Append new row 35, copying format from row 34 first ---
$ws.Range("A34:E34").Copy($ws.Range("A35:E35"))
$newRow = New-Object 'object[,]' 1,5
$newRow[0,0] = 33
$newRow[0,1] = 1485.514423856088
$newRow[0,2] = 5
$newRow[0,3] = 0.03083767517834564
$newRow[0,4] = 0.003172829395605959
$ws.Range("A35:E35").Value = $newRow

Write-Host "done"
